$d = $word.ActiveDocument

# The document starts with four header paragraphs:
#   1. Marlon Torres
#   2. 11/27/2013
#   3. Web Programming Fundamentals - Section 01
#   4. Activity: Problem Solving
# Change their font size from the default 12 pt to 14 pt.

for ($i = 1; $i -le 4; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Size = 14
    $p.Range.Font.SizeBi = 14
}
